{"js": "// Update the worksheet date and every \"NNN\u00d7N=\" multiplication prompt to the\n// new values from the target revision. Every old value is unique within the\n// document and no new value collides with a not-yet-processed old value, so\n// a simple sequential find/replace (matching full run text, case-sensitive)\n// is safe and deterministic.\nconst replacements = [\n  [\"2024-07-02 Tuesday\", \"2024-07-03 Wednesday\"],\n  [\"282\u00d78=\", \"966\u00d73=\"],\n  [\"690\u00d78=\", \"648\u00d79=\"],\n  [\"677\u00d72=\", \"965\u00d75=\"],\n  [\"747\u00d76=\", \"932\u00d78=\"],\n  [\"999\u00d72=\", \"672\u00d76=\"],\n  [\"894\u00d78=\", \"344\u00d73=\"],\n  [\"151\u00d79=\", \"381\u00d79=\"],\n  [\"864\u00d76=\", \"209\u00d78=\"],\n  [\"993\u00d77=\", \"142\u00d79=\"],\n  [\"974\u00d74=\", \"884\u00d76=\"],\n  [\"979\u00d75=\", \"474\u00d78=\"],\n  [\"409\u00d74=\", \"514\u00d74=\"],\n  [\"548\u00d79=\", \"273\u00d76=\"],\n  [\"256\u00d76=\", \"665\u00d78=\"],\n  [\"592\u00d79=\", \"880\u00d79=\"],\n  [\"581\u00d77=\", \"706\u00d79=\"],\n  [\"490\u00d74=\", \"717\u00d73=\"],\n  [\"968\u00d79=\", \"528\u00d73=\"],\n  [\"866\u00d72=\", \"119\u00d72=\"],\n  [\"375\u00d73=\", \"768\u00d76=\"],\n  [\"245\u00d76=\", \"236\u00d79=\"],\n  [\"663\u00d74=\", \"127\u00d76=\"],\n  [\"748\u00d79=\", \"785\u00d76=\"],\n  [\"364\u00d79=\", \"650\u00d78=\"],\n  [\"862\u00d74=\", \"768\u00d75=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"text\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date and every \"NNN\u00d7N=\" multiplication prompt to the\n# new values from the target revision. Every old value is unique within the\n# document and no new value collides with a not-yet-processed old value, so\n# a simple sequential Find/Replace (whole-document, case-sensitive, one\n# match each) is safe and deterministic.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"2024-07-02 Tuesday\", \"2024-07-03 Wednesday\"),\n  @(\"282\u00d78=\", \"966\u00d73=\"),\n  @(\"690\u00d78=\", \"648\u00d79=\"),\n  @(\"677\u00d72=\", \"965\u00d75=\"),\n  @(\"747\u00d76=\", \"932\u00d78=\"),\n  @(\"999\u00d72=\", \"672\u00d76=\"),\n  @(\"894\u00d78=\", \"344\u00d73=\"),\n  @(\"151\u00d79=\", \"381\u00d79=\"),\n  @(\"864\u00d76=\", \"209\u00d78=\"),\n  @(\"993\u00d77=\", \"142\u00d79=\"),\n  @(\"974\u00d74=\", \"884\u00d76=\"),\n  @(\"979\u00d75=\", \"474\u00d78=\"),\n  @(\"409\u00d74=\", \"514\u00d74=\"),\n  @(\"548\u00d79=\", \"273\u00d76=\"),\n  @(\"256\u00d76=\", \"665\u00d78=\"),\n  @(\"592\u00d79=\", \"880\u00d79=\"),\n  @(\"581\u00d77=\", \"706\u00d79=\"),\n  @(\"490\u00d74=\", \"717\u00d73=\"),\n  @(\"968\u00d79=\", \"528\u00d73=\"),\n  @(\"866\u00d72=\", \"119\u00d72=\"),\n  @(\"375\u00d73=\", \"768\u00d76=\"),\n  @(\"245\u00d76=\", \"236\u00d79=\"),\n  @(\"663\u00d74=\", \"127\u00d76=\"),\n  @(\"748\u00d79=\", \"785\u00d76=\"),\n  @(\"364\u00d79=\", \"650\u00d78=\"),\n  @(\"862\u00d74=\", \"768\u00d75=\")\n)\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $oldText\n  $find.Replacement.Text = $newText\n  $find.Forward = $true\n  $find.Wrap = 1\n  $find.Format = $false\n  $find.MatchCase = $true\n  $find.MatchWholeWord = $false\n  $find.MatchWildcards = $false\n\n  $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
